# Add 2022-Q1 data:
#  - new worksheet "2022-Q1" (fund holders detail) inserted right before "总计"
#  - "总计" (summary) worksheet gets a new first data row for 2022-Q1

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" worksheet, positioned right before "总计"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($total)
$newSheet.Name = "2022-Q1"

# Copy the cell formatting (styles/borders/fonts) from the "2021-Q4" sheet,
# which has the identical A1:H column layout, so the new sheet matches the
# existing look (bold/centered header row + index column). The template only
# has 6 data rows (A1:H6); extend the same row-6 formatting down to the
# extra rows 7:8 that this sheet needs.
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("A1:H6").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)
$template.Range("A6:H6").Copy()
$newSheet.Range("A7:H8").PasteSpecial(-4122)

# ---- header row ----
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---- data rows ----
# columns B,C,D,E,F,G are stored as text in this workbook's convention
# (fund codes keep leading zeros, percentages keep trailing zeros), so force
# a text number-format before assigning the values.
$dataRange = $newSheet.Range("B2:G8")
$dataRange.NumberFormat = "@"

$rows = @(
    @(0, "160425", "华安创业板两年定期开放混合", "5.11", "96.75", "3.00", "0.1533", 7),
    @(1, "003966", "中银润利灵活配置混合A",       "5.69", "25.96", "0.54", "0.0307", 10),
    @(2, "000059", "国联安中证医药100指数A",       "2.20", "91.27", "1.36", "0.0299", 7),
    @(3, "002614", "中银颐利灵活配置混合A",       "3.34", "36.07", "0.73", "0.0244", 9),
    @(4, "003967", "中银润利灵活配置混合C",       "3.85", "25.96", "0.54", "0.0208", 10),
    @(5, "002615", "中银颐利灵活配置混合C",       "2.43", "36.07", "0.73", "0.0177", 9),
    @(6, "006569", "国联安中证医药100指数C",       "0.34", "91.27", "1.36", "0.0046", 7)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Update the "总计" worksheet: add a 2022-Q1 row on top, push the rest
#    down and renumber the index column.
# ---------------------------------------------------------------------
# Re-fetch by name: the worksheet-collection insert above shifted "总计" to
# a new position, so the earlier $total reference no longer points at it.
$total = $wb.Worksheets.Item("总计")

# The table grows from 6 to 7 data rows; row 7 is brand new, so give its
# index cell (column A) the same style as the existing index column (copy
# the format already applied to A6) before the values get written below.
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)

$summaryRows = @(
    @(0, "2022-Q1", 7,  0.28),
    @(1, "2021-Q4", 5,  10.28),
    @(2, "2021-Q3", 13, 11.67),
    @(3, "2021-Q2", 17, 8.48),
    @(4, "2021-Q1", 14, 3.53),
    @(5, "2020-Q4", 40, 9.619999999999999)
)

$r = 2
foreach ($row in $summaryRows) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
